$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Anchor cell with known plain (style-0, text) formatting, used as a format-paste source
# so that forcing Text number-format on target cells can be cleanly undone afterwards.
$fmtSource = $ws.Range("B2")

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $fmtSource.Copy()
    $rng.PasteSpecial(-4122)
}

# --- Column D price updates (force-text to preserve literal formatting/leading zeros) ---
Set-TextValue "D2" "26.832.81"
Set-TextValue "D3" "1.870.96"
Set-TextValue "D5" "301.42"
Set-TextValue "D7" "0.5329"
Set-TextValue "D8" "0.3748"
Set-TextValue "D9" "0.07174"
Set-TextValue "D10" "21.49"
Set-TextValue "D11" "0.8864"
Set-TextValue "D13" "1.889.07"
Set-TextValue "D14" "93.18"
Set-TextValue "D15" "5.262"
Set-TextValue "D16" "1.002"
Set-TextValue "D17" "14.69"
Set-TextValue "D18" "0.000008535"
Set-TextValue "D19" "1.000"
Set-TextValue "D20" "26.884.27"
Set-TextValue "D21" "4.970"
Set-TextValue "D23" "6.395"
Set-TextValue "D24" "147.16"
Set-TextValue "D25" "2.256"
Set-TextValue "D26" "1.731"
Set-TextValue "D27" "18.02"
Set-TextValue "D28" "114.29"
Set-TextValue "D29" "4.733"
Set-TextValue "D30" "4.590"
Set-TextValue "D31" "0.09120"
Set-TextValue "D32" "0.8002"
Set-TextValue "D33" "0.04988"
Set-TextValue "D34" "1.175"
Set-TextValue "D35" "2.987"
Set-TextValue "D36" "0.6035"
Set-TextValue "D37" "2.589"
Set-TextValue "D38" "3.153"
Set-TextValue "D41" "6.592"
Set-TextValue "D42" "8.833"
Set-TextValue "D43" "115.82"
Set-TextValue "D44" "0.5155"
Set-TextValue "D45" "0.1495"
Set-TextValue "D46" "0.9997"
Set-TextValue "D47" "9.928"
Set-TextValue "D48" "1.627"
Set-TextValue "D49" "37.60"
Set-TextValue "D51" "62.11"

# --- Column B/C/E updates (already plain text, safe to assign directly) ---
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -6.03%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  +6.11%  "
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -6.27%  "
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("E51").Value = "  -2.87%  "

$excel.CutCopyMode = $false